$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Max N")

$ws.Range("B2").Value = 22
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 17

$ws.Activate()
$ws.Range("F2").Select()
